{"js": "// Apply the benchmark-stats fixups to the single-column results table.\n// Each table row holds one metric value in its only cell; a handful of\n// rows are rewritten outright (single value -> single value) and the\n// last three rows collapse a tab-separated 10-value line down to just\n// its first value.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index (0-based) -> new cell text\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"211\",\n  4: \"0.00002\",\n  5: \"0.00019\",\n  6: \"0.00004\",\n  7: \"0.00001\",\n  11: \"0.00786\",\n  43: \"100\",\n  44: \"0.01\",\n  45: \"186\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-stats fixups to the single-column results table.\n# Each table row holds one metric value in its only cell; a handful of\n# rows are rewritten outright (single value -> single value) and the\n# last three rows collapse a tab-separated 10-value line down to just\n# its first value.\n$doc = $word.ActiveDocument\n$tbl = $doc.Tables.Item(1)\n\n# Word row index is 1-based; map (1-based row) -> new cell text.\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"211\"\n    5  = \"0.00002\"\n    6  = \"0.00019\"\n    7  = \"0.00004\"\n    8  = \"0.00001\"\n    12 = \"0.00786\"\n    44 = \"100\"\n    45 = \"0.01\"\n    46 = \"186\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $tbl.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
